# Update cryptocurrency price/volume values on Sheet1, matching the
# Dec 28 2022 GitHub Actions symbol-list refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "243.65"
    "D4" = "5.249"
    "D5" = "0.05814"
    "D6" = "6.487"
    "D7" = "3.331"
    "D8" = "0.8083"
    "D9" = "0.8747"
    "D11" = "0.07276"
    "D12" = "0.03056"
    "D13" = "0.03055"
    "D14" = "0.09317"
    "D16" = "0.001539"
    "D17" = "0.04704"
    "D18" = "0.0006008"
    "E18" = "17OneONEWorstin24h"
    "D19" = "0.006112"
    "D20" = "0.001266"
    "D21" = "0.004591"
    "D22" = "0.00008696"
    "D24" = "2.147"
    "D40" = "0.03794"
    "D41" = "0.006287"
    "D42" = "0.1054"
    "D43" = "0.002698"
    "E43" = "42CEJICEJI"
    "D44" = "0.007640"
    "D45" = "0.00005465"
    "D47" = "0.5498"
    "D48" = "0.006980"
    "D49" = "0.00002099"
    "D50" = "0.0001999"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so the numeric-looking strings keep their
    # exact printed form (trailing zeros, precision) instead of being
    # auto-coerced to a Double by Excel's smart-entry parser.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Drop back to the default "Normal" style so no stray numFmt/style
    # gets attached to the cell (matches original formatting).
    $cell.Style = "Normal"
}
